# This script applies a set of spelling/typo corrections to the text
# contained in several cells of Sheet1, and updates a couple of
# view-related settings (window size, top-left visible cell, and the
# current selection) to match the author's saved workbook state.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix typos in cell text (same cell positions, corrected spelling) ---

$ws.Range("F1").Value = "Recommendation"

$ws.Range("C2").Value = "As per observation made the web application OWASP Juice Shop is vulnerable to SQL injection, Where the attacker will enter a malicious code instead of correct data in login credentials through which the attacker will gain access as a authenticated user"

$ws.Range("F2").Value = "1. Use safe API or adopt the use of Object Relational mapping Tools.                                                                                               2.Use Server side validation and prevent use of userID's and passwords using special characters."

$ws.Range("C7").Value = "when the packets are in intercepted  using burp suite the data in packets are sent in plain text."

$ws.Range("D7").Value = "The most basic flaw is simply not encrypting data manual attack is generally required. Previously retrieved password databases could be brute forced by Graphics Processing Units (GPUs)."

$ws.Range("D10").Value = "These flaws can be used to extract data, execute a remote request from the server, scan internal systems, perform a denial-of-service attack, as well as execute other attacks. The business impact depends in the protection reeds of all affected application and data."

$ws.Range("F10").Value = "1.Avoid serialization of data.                                                             2.Implement server ide validation, filtering or sanitization."

$ws.Range("C11").Value = "As per the observation the script can be executed in the URL."

$ws.Range("F11").Value = "1.Apply context sensitive encoding when modifying the browser document on the client side"

# --- Update the view state: scroll position and current selection ---
# (Best effort: the saved view shows the A10:F11 block selected, scrolled so
# row 5 is at the top of the window.)

$window = $excel.ActiveWindow
$window.ScrollRow = 5
$window.ScrollColumn = 1

$ws.Range("A10:F11").Select()
